$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update LakePerimeter (B3) and LakeArea (B5) values ---
$ws.Range("B3").Value2 = 25900
$ws.Range("B5").Value2 = 16079000

# --- LakeVolume (B6) becomes a computed formula (=LakeDepth*LakeArea) ---
$ws.Range("B6").Formula = "=B4*B5"

# Mark the recalculated cell with an explicit (no-op) fill application,
# matching the extra cellXfs record introduced upstream.
$ws.Range("B6").Interior.ColorIndex = -4142

# --- Column B width: new column was given a fitted width of 10 ---
$ws.Columns(2).ColumnWidth = 9.166666666666666

# --- Row 19 ("AerialLoad") was cut and moved up into the previously
#     empty row 18, leaving row 19 blank again ---
$ws.Range("A19:C19").Cut($ws.Range("A18:C18")) | Out-Null

# --- Update the active selection on the sheet ---
$ws.Range("B16").Select() | Out-Null
